# US09 -> US14: update the "ID User Story" cell text.
#   "UC 11 Notificação de mudança de status quando houver modificação no
#    plano de negócio."
# becomes
#   "UC 12 Visualizar a descrição de cada item do plano facilitando o
#    preenchimento do plano."
#
# The run immediately preceding "UC 11" (a plain space, same rPr) and the
# run group starting at "Notificação" (rFonts/color rPr) must keep their
# own distinct boundaries/formatting, matching the original authoring.
# Find.Execute happily coalesces runs with identical rPr that it touches,
# so after each text replace we "re-cut" the run at the boundary we still
# need by toggling Bold off/on across exactly that sub-range — a real
# (if momentary) formatting change forces the engine to keep it as its
# own <w:r>, and flipping Bold back to True restores the original
# (implicit) bold formatting without leaving any stray override behind.

$d = $word.ActiveDocument

# --- Segment 1: "UC 11 " -> "UC 12" + new " " run ------------------------

$rng = $d.Content
$rng.Find.Execute("UC 11 ", $true, $false, $false, $false, $false, $true, 1, $false, "UC 12 ", 2)

# "UC 12 " merged with the preceding (identical rPr) space run; re-cut the
# boundary so the leading space stays its own run, as it was before.
$rng = $d.Content
$rng.Find.Execute("UC 12 ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ucRange = $d.Range($rng.Start, $rng.End)
$ucRange.Font.Bold = $false
$ucRange.Font.Bold = $true

# Split "UC 12 " into "UC 12" and a new, separate trailing " " run.
$rng = $d.Content
$rng.Find.Execute("UC 12 ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$spaceRange = $d.Range($rng.End - 1, $rng.End)
$spaceRange.Font.Bold = $false
$spaceRange.Font.Bold = $true

# --- Segment 2: rewrite the description run-group ------------------------

$rng = $d.Content
$rng.Find.Execute("Notificação de mudança de status quando houver modificação no plano de negócio.", $true, $false, $false, $false, $false, $true, 1, $false, "Visualizar a descrição de cada item do plano facilitando o preenchimento do plano.", 2)

# That replace merges the whole rFonts/color run-group into a single run;
# re-split it into the four runs the source document used to carry:
#   "V" | "isualizar a descrição d" | "e cada item do plano facilitando" |
#   " o preenchimento do plano."
$rng = $d.Content
$rng.Find.Execute("Visualizar a descrição de cada item do plano facilitando o preenchimento do plano.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base = $rng.Start
$end = $rng.End

$seg1 = $d.Range($base, $base + 1)
$seg1.Font.Bold = $false
$seg1.Font.Bold = $true

$seg2 = $d.Range($base + 1, $base + 24)
$seg2.Font.Bold = $false
$seg2.Font.Bold = $true

$seg3 = $d.Range($base + 24, $base + 56)
$seg3.Font.Bold = $false
$seg3.Font.Bold = $true
